# Rename the header row (A1:U1) so the "old"/"new" suffixes are replaced with
# the actual format-version identifiers the two compared AHB files came from
# (FV2310 = "old" / left-hand side, FV2404 = "new" / right-hand side), then
# turn the sheet's data range into a real Excel Table ("Table1") and freeze
# the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headersBefore = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

$headersAfter = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

# Columns A-J: "<name>_old" -> "<name>_FV2310"
for ($i = 0; $i -lt $headersBefore.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersBefore[$i]
}

# Column K stays "diff"
$ws.Cells.Item(1, 11).Value = "diff"

# Columns L-U: "<name>_new" -> "<name>_FV2404"
for ($i = 0; $i -lt $headersAfter.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersAfter[$i]
}

# Turn the used range into an Excel Table, matching the exported structure
# (ref A1:U73, 21 columns, one per header above).
$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header renaming, table creation and freeze pane applied."
